# Update coordinates for Whiteley map
$wb = $excel.ActiveWorkbook

# --- "d coordinates" sheet: insert an "ID" column (B) between Point and Latitude ---
$ws2 = $wb.Worksheets.Item("d coordinates")

$ws2.Range("B1").EntireColumn.Insert()

$ws2.Range("B1").Value = "ID"

$ids = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws2.Range("B$row").Value = $ids[$i]
}

# --- Add the new "whiteley coordinates" sheet after "d coordinates" ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $afterSheet)
$ws3.Name = "whiteley coordinates"

$ws3.Range("A1").Value = "Point"
$ws3.Range("B1").Value = "ID"
$ws3.Range("C1").Value = "Latitude"
$ws3.Range("D1").Value = "Longitude"
$ws3.Range("A1:B1").NumberFormat = "@"

$names = @(
    "red_start","red_midpoint","red_end",
    "orange_branchpoint","orange_end_1","orange_end_2",
    "yellow_start","yellow_end",
    "green_start","green_branchpoint","green_end_1","green_end_2",
    "pink_start","pink_branchpoint_1","pink_end_1","pink_branchpoint_2","pink_end_2","pink_end_3",
    "blue_start","blue_branchpoint_1","blue_end_1","blue_branchpoint_2","blue_end_2","blue_end_3"
)
$idvals = @(1,2,3,4,5,6,7,9,8,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws3.Range("A$row").Value = $names[$i]
    $ws3.Range("B$row").Value = $idvals[$i]
}

# --- View / zoom tweaks ---
$ws1 = $wb.Worksheets.Item("c coordinates")
$ws1.Select()
$excel.ActiveWindow.Zoom = 115
$ws1.Range("F4").Select()

$ws2.Select()
$excel.ActiveWindow.Zoom = 115
$ws2.Range("A1:D1").Select()

$ws3.Select()
$ws3.Range("H9").Select()
